$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.233"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05773"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.435"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.233"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = "'0.8797"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1384"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07103"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03158"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03042"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09318"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.821"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001522"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04725"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0006016"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006205"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001260"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004056"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00008719"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.542"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.151"
$ws.Range("D24").Style = "Normal"
$ws.Range("D28").Value = "'0.0002335"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.03733"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006284"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1044"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002482"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.007144"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005337"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.5364"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002558"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002105"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002005"
$ws.Range("D50").Style = "Normal"
